# Calculated total sal and annual sal of the employees
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Headers
$ws.Range("I1").Value = "Total"
$ws.Range("J1").Value = "Annual sal"

# Row 2 formulas entered individually (become plain, non-shared formulas)
$ws.Range("I2").Formula = "=E2+H2"
$ws.Range("J2").Formula = "=I2*12"

# Rows 3-22 filled as a block (become a shared formula group)
$ws.Range("I3:I22").Formula = "=E3+H3"
$ws.Range("J3:J22").Formula = "=I3*12"

# Selection as in the diff (activeCell J2, sqref J2:J22)
$ws.Range("J2:J22").Select()
